# StructureDefinition-ror-organization-reopening-date.xlsx edit
#
# 1) Metadata sheet: bump the "Date" value (row 8, column B).
# 2) Elements sheet: the two "Mapping: ..." columns (AK = col 37,
#    AL = col 38) trade places - header text, per-row data, and
#    (best-effort) column width all swap together, while every other
#    column / row stays untouched.

$wb = $excel.ActiveWorkbook

# --- 1) Metadata!B8 : Date -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- 2) Elements!AK<->AL : swap the two mapping columns --------------------
$elements = $wb.Worksheets.Item("Elements")

# Swap header text (row 1)
$akHeader = $elements.Range("AK1").Value2
$alHeader = $elements.Range("AL1").Value2
$elements.Range("AK1").Value = $alHeader
$elements.Range("AL1").Value = $akHeader

# Swap the data rows (2 through 6 - the full used range below the header)
for ($row = 2; $row -le 6; $row++) {
    $akCell = $elements.Cells.Item($row, 37)
    $alCell = $elements.Cells.Item($row, 38)

    $akValue = $akCell.Value2
    $alValue = $alCell.Value2

    if ($null -eq $alValue) {
        $akCell.Value = ""
    } else {
        $akCell.Value = $alValue
    }

    if ($null -eq $akValue) {
        $alCell.Value = ""
    } else {
        $alCell.Value = $akValue
    }
}

# Swap the (best-fit) column widths so AK/AL keep matching their new,
# swapped header text. (Target stored widths: 70.94140625 / 24.98046875 -
# the engine quantises ColumnWidth writes, so these inputs are chosen to
# land on the closest representable value.)
$elements.Columns.Item(37).ColumnWidth = 70.1
$elements.Columns.Item(38).ColumnWidth = 24.1
